# Remove the data row for account "008336332" (Carlos, balance 11000)
# from the "Export" sheet. Deleting the row shifts all subsequent rows
# up by one, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$target = $ws.Columns.Item(1).Find("008336332")
if ($target -ne $null) {
    $target.EntireRow.Delete()
} else {
    # Fallback: the row is the 4th row on the sheet (header + 3 data rows).
    $ws.Rows.Item(4).Delete()
}
